$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width (closest achievable to 23.08984375 given engine's internal
# pixel-rounding of ColumnWidth)
$ws.Columns.Item(1).ColumnWidth = 23

# --- New content: "Display" section parts ---
# The order in which NEW unique strings are introduced below matters: it
# determines the order they are appended to xl/sharedStrings.xml (indices
# 9..20), matching the target workbook.

# index 9
$ws.Range("A19").Value = "Display:"

# index 10
$ws.Range("A21").Value = "Multi Meters"

# index 11 (value + hyperlink, same URL text)
$ws.Range("B21").Value = "https://www.amazon.com/McIgIcM-Digital-Voltmeter-Ammeter-10ADetector/dp/B06XR2XKNT/ref=sr_1_7?keywords=digital+current+meter&qid=1561517841&s=industrial&sr=1-7"
$ws.Range("B21").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B21"), "https://www.amazon.com/McIgIcM-Digital-Voltmeter-Ammeter-10ADetector/dp/B06XR2XKNT/ref=sr_1_7?keywords=digital+current+meter&qid=1561517841&s=industrial&sr=1-7")
$ws.Range("B21").Style = "Hyperlink"

# index 12 (value + hyperlink, same URL text) -- set on B23 before A23 so
# that the shared-string order matches the target exactly
$ws.Range("B23").Value = "https://www.digikey.com/product-detail/en/adafruit-industries-llc/705/1528-1149-ND/5353609"
$ws.Range("B23").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B23"), "https://www.digikey.com/product-detail/en/adafruit-industries-llc/705/1528-1149-ND/5353609")
$ws.Range("B23").Style = "Hyperlink"

# index 13
$ws.Range("A23").Value = "Volt Meters"

# index 14 (value + hyperlink, same URL text) -- set on B25 before A25
$ws.Range("B25").Value = "http://easycircuit012.blogspot.com/2012/12/digital-volt-and-ampere-meter-circuit.html"
$ws.Range("B25").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B25"), "http://easycircuit012.blogspot.com/2012/12/digital-volt-and-ampere-meter-circuit.html")
$ws.Range("B25").Style = "Hyperlink"

# index 15
$ws.Range("A25").Value = "Design it on your own"

# index 16 (value + hyperlink, same URL text)
$ws.Range("B27").Value = "https://www.digikey.com/product-detail/en/lumex-opto-components-inc/LDT-A512RI/67-1424-ND/252626"
$ws.Range("B27").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B27"), "https://www.digikey.com/product-detail/en/lumex-opto-components-inc/LDT-A512RI/67-1424-ND/252626")
$ws.Range("B27").Style = "Hyperlink"

# index 17 (value + hyperlink, same URL text)
$ws.Range("B28").Value = "https://www.digikey.com/product-detail/en/lumex-opto-components-inc/LDT-A514RI/67-1423-ND/252628"
$ws.Range("B28").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B28"), "https://www.digikey.com/product-detail/en/lumex-opto-components-inc/LDT-A514RI/67-1423-ND/252628")
$ws.Range("B28").Style = "Hyperlink"

# index 18 (value + hyperlink, same URL text)
$ws.Range("B30").Value = "http://www.electronics-diy.com/70v_pic_voltmeter_amperemeter.php"
$ws.Range("B30").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B30"), "http://www.electronics-diy.com/70v_pic_voltmeter_amperemeter.php")
$ws.Range("B30").Style = "Hyperlink"

# index 19 (value + hyperlink, same URL text)
$ws.Range("B31").Value = "https://www.adafruit.com/product/399?gclid=Cj0KCQjwjMfoBRDDARIsAMUjNZpkoKOfF8LGOWSxf4VLWIPIVsbuiiDFCeu1C8yUAKVEhcIb9canXt8aAnndEALw_wcB"
$ws.Range("B31").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B31"), "https://www.adafruit.com/product/399?gclid=Cj0KCQjwjMfoBRDDARIsAMUjNZpkoKOfF8LGOWSxf4VLWIPIVsbuiiDFCeu1C8yUAKVEhcIb9canXt8aAnndEALw_wcB")
$ws.Range("B31").Style = "Hyperlink"

# index 20 (value + hyperlink, same URL text)
$ws.Range("B32").Value = "https://www.digikey.com/product-detail/en/PIC16F876A-I%2fSO/PIC16F876A-I%2fSO-ND/446139/?itemSeq=296474682"
$ws.Range("B32").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B32"), "https://www.digikey.com/product-detail/en/PIC16F876A-I%2fSO/PIC16F876A-I%2fSO-ND/446139/?itemSeq=296474682")
$ws.Range("B32").Style = "Hyperlink"

# Final selection matches the saved view in the target workbook
$ws.Range("B32").Select()
